# CIV-11205: replace the <<courtName>> merge field in the "This order is
# made by ..." recital with <<siteName>> - <<address>> - <<postcode>>.
#
# The document contains another, unrelated "<<courtName>>" occurrence in
# the heading table, so everything below is scoped to a small Range around
# the unique phrase "This order is made by" to avoid touching that one.

$d = $word.ActiveDocument

# 1. Locate the target paragraph via a phrase that only appears once.
$full = $d.Content
$found = $full.Find.Execute("This order is made by", $true, $false, $false, `
    $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find target paragraph"
}
$paraStart = $full.Start
$para = $d.Range($paraStart, $paraStart + 200)

# 2. Within that paragraph, find "courtName" (search only - no in-place
#    Find/Replace, which in this runtime can mutate the first match of the
#    literal string anywhere in the document rather than the scoped one).
$search = $para.Duplicate
$found2 = $search.Find.Execute("courtName", $false, $false, $false, $false, `
    $false, $false, 1, $false, "", 0)
if (-not $found2) {
    throw "Could not find courtName field in target paragraph"
}

# Rename the field in place - this keeps the existing run (and the
# surrounding spell-check proofErr markers) untouched, only the text
# changes.
$courtNameRange = $d.Range($search.Start, $search.End)
$courtNameRange.Text = "siteName"

# 3. Find the closing ">>." that follows (now shifted because "courtName"
#    -> "siteName" is one character shorter) and insert the new address /
#    postcode fields just before the final full stop.
$search2 = $para.Duplicate
$found3 = $search2.Find.Execute(">>.", $false, $false, $false, $false, `
    $false, $false, 1, $false, "", 0)
if (-not $found3) {
    throw "Could not find end of courtName/siteName field"
}

$insertPoint = $d.Range($search2.Start + 2, $search2.Start + 2)
$insertPoint.InsertAfter(" - <<address>> - <<postcode>>")
